$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Model")

# Add the new Q4 FY24 quarter column (column I) header
$ws.Range("I1").Value = "Q4 FY24"

# Report Date for the new quarter (I2) - give it the same date format as the
# neighbouring date cells by copying formats only (keeps the shared style
# instead of minting a brand new custom number format).
$ws.Range("I2").Value = 45654
$ws.Range("H2").Copy() | Out-Null
$ws.Range("I2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# New quarter figures: Data Center, Client, Gaming, Embedded
$ws.Range("I3").Value = 3859
$ws.Range("I4").Value = 2313
$ws.Range("I5").Value = 563
$ws.Range("I6").Value = 923

# The new date value in column I is now as wide as the other "longer" date
# columns (e.g. J), so Excel's column auto-fit widens it to match.
$ws.Columns.Item(9).ColumnWidth = 9.25

# Update the active selection to match the authored state
$ws.Range("I6").Select()
